# factor(one_file_one_tab): end factorisation of InsertController methods
#
# Adds two new worksheets ("time_min" and "time_min_expected") at the end of
# the workbook, each holding the small "temps utilise" sample table, the
# second one with an extra computed column (time expressed in minutes).
# Each new sheet also gets a little bar chart (as the other data sheets in
# this workbook do), with the "time_min_expected" sheet getting two of them.

$wb = $excel.ActiveWorkbook

function Fill-TimeSheet($ws, [bool]$withMinutes) {
    $ws.Range("A1").Value = "Nom"
    $ws.Range("B1").Value = "Prénom"
    $ws.Range("C1").Value = "Adresse de courriel"
    $ws.Range("D1").Value = "État"
    $ws.Range("E1").Value = "Temps utilisé"

    $ws.Range("A2").Value = "Abdulhoussen"
    $ws.Range("B2").Value = "Houzefa"
    $ws.Range("C2").Value = "houzefa.abdulhoussen@universite-paris-saclay.fr"
    $ws.Range("D2").Value = "Terminé"
    $ws.Range("E2").Value = "5 min 49 s"

    $ws.Range("A3").Value = "Abdel Moneim"
    $ws.Range("B3").Value = "Yasmine"
    $ws.Range("C3").Value = "yasmine.abdel-moneim@universite-paris-saclay.fr"
    $ws.Range("D3").Value = "Terminé"
    $ws.Range("E3").Value = "7 min"

    $ws.Range("A4").Value = "Abbas"
    $ws.Range("B4").Value = "Zina"
    $ws.Range("C4").Value = "zina.abbas@universite-paris-saclay.fr"
    $ws.Range("D4").Value = "Terminé"
    $ws.Range("E4").Value = "1 jour 5 heures"

    $ws.Range("A5").Value = "Abdallah"
    $ws.Range("B5").Value = "Aboubaker"
    $ws.Range("C5").Value = "aboubaker.abdallah@universite-paris-saclay.fr"
    $ws.Range("D5").Value = "Terminé"
    $ws.Range("E5").Value = "2 jours 10 min"

    if ($withMinutes) {
        $ws.Range("F2").Value = "5,82"
        $ws.Range("F3").Value = "7,0"
        $ws.Range("F4").Value = "1740,0"
        $ws.Range("F5").Value = "2890,0"
    }
}

function Add-TimeChart($ws, $sheetName) {
    $chartObj = $ws.ChartObjects().Add(10, 10, 430, 300)
    $chart = $chartObj.Chart
    $chart.ChartType = -4100
    $chart.HasTitle = $true
    $chart.ChartTitle.Text = "taux réussite ( %) par question"

    $series = $chart.SeriesCollection()
    $newSeries = $series.NewSeries()
    $newSeries.Name = "=" + $sheetName + "!`$E`$1"
    $newSeries.XValues = $ws.Range("A2:A5")
    $newSeries.Values = $ws.Range("E2:E5")
}

# --- add the two new sheets at the end of the workbook ---------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$timeMin = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$timeMin.Name = "time_min"
Fill-TimeSheet $timeMin $false
Add-TimeChart $timeMin "time_min"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$timeMinExpected = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$timeMinExpected.Name = "time_min_expected"
Fill-TimeSheet $timeMinExpected $true
Add-TimeChart $timeMinExpected "time_min_expected"
Add-TimeChart $timeMinExpected "time_min_expected"

# "time_min" ends up being the active tab, matching the recorded edit.
$timeMin.Select()
